$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.201571333333333
$ws.Range("H2").Value = 18.604714
$ws.Range("I2").Value = 0.05221490529364391
$ws.Range("J2").Value = 0.07406232529850043
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.68421466666667
$ws.Range("N2").Value = 32.052644
$ws.Range("O2").Value = 0.04705285980693976
$ws.Range("P2").Value = 0.04892736897547583
$ws.Range("Q2").Value = 66.25891939597956
$ws.Range("R2").Value = 596.3302745638159
$ws.Range("S2").Value = 0.002456860618614463
$ws.Range("T2").Value = 0.003623674717061449
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.201571333333333
$ws.Range("H3").Value = 18.604714
$ws.Range("I3").Value = 0.05221490529364391
$ws.Range("J3").Value = 0.07406232529850043
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 70.36235166666667
$ws.Range("N3").Value = 211.087055
$ws.Range("O3").Value = 0.3098730203341347
$ws.Range("P3").Value = 0.3222178559101571
$ws.Range("Q3").Value = 436.3571430419188
$ws.Range("R3").Value = 3927.21428737727
$ws.Range("S3").Value = 0.01617999040980224
$ws.Range("T3").Value = 0.02386420366140339
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.201571333333333
$ws.Range("H4").Value = 18.604714
$ws.Range("I4").Value = 0.05221490529364391
$ws.Range("J4").Value = 0.07406232529850043
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 59.09107466666666
$ws.Range("N4").Value = 177.273224
$ws.Range("O4").Value = 0.2602347611759026
$ws.Range("P4").Value = 0.2706020894912812
$ws.Range("Q4").Value = 366.4575147086595
$ws.Range("R4").Value = 3298.117632377935
$ws.Range("S4").Value = 0.0135881334089138
$ws.Range("T4").Value = 0.0200414199783572
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.201571333333333
$ws.Range("H5").Value = 18.604714
$ws.Range("I5").Value = 0.05221490529364391
$ws.Range("J5").Value = 0.07406232529850043
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 60.83231733333333
$ws.Range("N5").Value = 182.496952
$ws.Range("O5").Value = 0.2679031251727568
$ws.Range("P5").Value = 0.2785759485989269
$ws.Range("Q5").Value = 377.2559553146364
$ws.Range("R5").Value = 3395.303597831728
$ws.Range("S5").Value = 0.01398853630876673
$ws.Range("T5").Value = 0.02063198252547206
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.201571333333333
$ws.Range("H6").Value = 18.604714
$ws.Range("I6").Value = 0.05221490529364391
$ws.Range("J6").Value = 0.07406232529850043
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 26.0983795
$ws.Range("N6").Value = 52.196759
$ws.Range("O6").Value = 0.1149362335102661
$ws.Range("P6").Value = 0.07967673702415903
$ws.Range("Q6").Value = 161.8509621536543
$ws.Range("R6").Value = 971.1057729219259
$ws.Range("S6").Value = 0.006001384547546688
$ws.Range("T6").Value = 0.00590104441620634
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.461641333333333
$ws.Range("H7").Value = 22.384924
$ws.Range("I7").Value = 0.06282422221945559
$ws.Range("J7").Value = 0.0891107233935555
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.68421466666667
$ws.Range("N7").Value = 32.052644
$ws.Range("O7").Value = 0.04705285980693976
$ws.Range("P7").Value = 0.04892736897547583
$ws.Range("Q7").Value = 79.72177777100622
$ws.Range("R7").Value = 717.495999939056
$ws.Range("S7").Value = 0.002956059320572074
$ws.Range("T7").Value = 0.004359953243148055
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.461641333333333
$ws.Range("H8").Value = 22.384924
$ws.Range("I8").Value = 0.06282422221945559
$ws.Range("J8").Value = 0.0891107233935555
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 70.36235166666667
$ws.Range("N8").Value = 211.087055
$ws.Range("O8").Value = 0.3098730203341347
$ws.Range("P8").Value = 0.3222178559101571
$ws.Range("Q8").Value = 525.0186315065355
$ws.Range("R8").Value = 4725.16768355882
$ws.Range("S8").Value = 0.01946753148928556
$ws.Range("T8").Value = 0.02871306623047453
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.461641333333333
$ws.Range("H9").Value = 22.384924
$ws.Range("I9").Value = 0.06282422221945559
$ws.Range("J9").Value = 0.0891107233935555
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 59.09107466666666
$ws.Range("N9").Value = 177.273224
$ws.Range("O9").Value = 0.2602347611759026
$ws.Range("P9").Value = 0.2706020894912812
$ws.Range("Q9").Value = 440.9164051638862
$ws.Range("R9").Value = 3968.247646474976
$ws.Range("S9").Value = 0.01634904646534186
$ws.Range("T9").Value = 0.02411354794637571
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.461641333333333
$ws.Range("H10").Value = 22.384924
$ws.Range("I10").Value = 0.06282422221945559
$ws.Range("J10").Value = 0.0891107233935555
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 60.83231733333333
$ws.Range("N10").Value = 182.496952
$ws.Range("O10").Value = 0.2679031251727568
$ws.Range("P10").Value = 0.2785759485989269
$ws.Range("Q10").Value = 453.9089334168497
$ws.Range("R10").Value = 4085.180400751648
$ws.Range("S10").Value = 0.0168308054691399
$ws.Range("T10").Value = 0.02482410429969631
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 7.461641333333333
$ws.Range("H11").Value = 22.384924
$ws.Range("I11").Value = 0.06282422221945559
$ws.Range("J11").Value = 0.0891107233935555
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 26.0983795
$ws.Range("N11").Value = 52.196759
$ws.Range("O11").Value = 0.1149362335102661
$ws.Range("P11").Value = 0.07967673702415903
$ws.Range("Q11").Value = 194.7367472102193
$ws.Range("R11").Value = 1168.420483261316
$ws.Range("S11").Value = 0.007220779475116199
$ws.Range("T11").Value = 0.007100051673860898
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 105.106922
$ws.Range("H12").Value = 210.213844
$ws.Range("I12").Value = 0.8849608724869005
$ws.Range("J12").Value = 0.836826951307944
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.68421466666667
$ws.Range("N12").Value = 32.052644
$ws.Range("O12").Value = 0.04705285980693976
$ws.Range("P12").Value = 0.04892736897547583
$ws.Range("Q12").Value = 1122.984917600589
$ws.Range("R12").Value = 6737.909505603536
$ws.Range("S12").Value = 0.04163993986775322
$ws.Range("T12").Value = 0.04094374101526633
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 105.106922
$ws.Range("H13").Value = 210.213844
$ws.Range("I13").Value = 0.8849608724869005
$ws.Range("J13").Value = 0.836826951307944
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 70.36235166666667
$ws.Range("N13").Value = 211.087055
$ws.Range("O13").Value = 0.3098730203341347
$ws.Range("P13").Value = 0.3222178559101571
$ws.Range("Q13").Value = 7395.570208364903
$ws.Range("R13").Value = 44373.42125018941
$ws.Range("S13").Value = 0.2742254984350469
$ws.Range("T13").Value = 0.2696405860182791
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 105.106922
$ws.Range("H14").Value = 210.213844
$ws.Range("I14").Value = 0.8849608724869005
$ws.Range("J14").Value = 0.836826951307944
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 59.09107466666666
$ws.Range("N14").Value = 177.273224
$ws.Range("O14").Value = 0.2602347611759026
$ws.Range("P14").Value = 0.2706020894912812
$ws.Range("Q14").Value = 6210.880975885509
$ws.Range("R14").Value = 37265.28585531306
$ws.Range("S14").Value = 0.230297581301647
$ws.Range("T14").Value = 0.2264471215665483
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 105.106922
$ws.Range("H15").Value = 210.213844
$ws.Range("I15").Value = 0.8849608724869005
$ws.Range("J15").Value = 0.836826951307944
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 60.83231733333333
$ws.Range("N15").Value = 182.496952
$ws.Range("O15").Value = 0.2679031251727568
$ws.Range("P15").Value = 0.2785759485989269
$ws.Range("Q15").Value = 6393.897633033914
$ws.Range("R15").Value = 38363.38579820348
$ws.Range("S15").Value = 0.2370837833948502
$ws.Range("T15").Value = 0.2331198617737585
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 105.106922
$ws.Range("H16").Value = 210.213844
$ws.Range("I16").Value = 0.8849608724869005
$ws.Range("J16").Value = 0.836826951307944
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 26.0983795
$ws.Range("N16").Value = 52.196759
$ws.Range("O16").Value = 0.1149362335102661
$ws.Range("P16").Value = 0.07967673702415903
$ws.Range("Q16").Value = 2743.120338432899
$ws.Range("R16").Value = 10972.48135373159
$ws.Range("S16").Value = 0.1017140694876033
$ws.Range("T16").Value = 0.0666756409340918
